$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RAF")

# Update Regional Availability Factor values per calibration
$ws.Range("B2").Value = 0.6
$ws.Range("B3").Value = 0.9
$ws.Range("B4").Value = 0.9
$ws.Range("B12").Value = 0.9
$ws.Range("B13").Value = 0.9
$ws.Range("B14").Value = 0.6
$ws.Range("B15").Value = 0.9
$ws.Range("B16").Value = 0.9
$ws.Range("B17").Value = 0.9
$ws.Range("B18").Value = 0.9
$ws.Range("B19").Value = 0.9
$ws.Range("B20").Value = 0.9
$ws.Range("B21").Value = 0.9
$ws.Range("B22").Value = 0.9
$ws.Range("B23").Value = 0.9
$ws.Range("B24").Value = 0.9

# Update the selection on the RAF sheet to B2:B24 with active cell B2,
# then restore the originally active sheet ("About") so the active tab
# is unchanged.
$ws.Select()
$ws.Range("B2:B24").Select()
$wb.Worksheets.Item("About").Select()
